{"js": "const body = context.document.body;\n\n// Step 1: the trailing \" 2\" (space + \"2\") in the first paragraph becomes \" 1\".\n// Searching for the two characters together and replacing them in one shot\n// merges the \" \" run and the \"2\" run into a single \" 1\" run, matching the\n// target OOXML (the \"demo\" run right before it is left untouched).\nconst numResults = body.search(\" 2\", { matchCase: true });\nnumResults.load(\"items\");\nawait context.sync();\n\nif (numResults.items.length > 0) {\n  numResults.items[0].insertText(\" 1\", \"Replace\");\n  await context.sync();\n}\n\n// Step 2: the second paragraph (\"Th\u00eam m\u1ed9t d\u00f2ng\") is merged back into the\n// first paragraph - i.e. the paragraph break between them is removed, while\n// the bookmark that lived at the end of paragraph 2 is preserved.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 1) {\n  const firstPara = paragraphs.items[0];\n  const secondPara = paragraphs.items[1];\n\n  // The range spanning from the end of paragraph 1 to the start of\n  // paragraph 2 is exactly the paragraph mark between them; deleting it\n  // joins the two paragraphs while keeping everything else (including the\n  // bookmark at the tail of paragraph 2) intact.\n  const paraMark = firstPara.getRange(\"End\").expandTo(secondPara.getRange(\"Start\"));\n  paraMark.delete();\n  await context.sync();\n\n  // Step 3: drop the leftover text that used to be the second paragraph's\n  // content (\"Th\u00eam m\u1ed9t d\u00f2ng\"), leaving the bookmark in place at the end.\n  const textResults = body.search(\"Th\u00eam m\u1ed9t d\u00f2ng\", { matchCase: true });\n  textResults.load(\"items\");\n  await context.sync();\n\n  if (textResults.items.length > 0) {\n    textResults.items[0].insertText(\"\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Step 1: change the trailing \"2\" in the first paragraph to \"1\" (the space\n# that already precedes it stays, and the write naturally merges the \" \"\n# run with the replaced run into a single \" 1\" run - the preceding \"demo\"\n# run is left untouched).\n$p1 = $d.Paragraphs(1)\n$findRange = $p1.Range.Duplicate\n$find = $findRange.Find\n$find.ClearFormatting()\n$find.Text = \"2\"\n$found = $find.Execute()\nif ($found) {\n    $findRange.Text = \"1\"\n}\n\n# Step 2: merge the second paragraph (\"Th\u00eam m\u1ed9t d\u00f2ng\") back into the first\n# one, removing its visible text but keeping the bookmark that sits at the\n# end of it.\nif ($d.Paragraphs.Count -gt 1) {\n    $p2 = $d.Paragraphs(2)\n\n    # Clear the second paragraph's text content (everything except its\n    # trailing paragraph mark, so the bookmark right after stays put).\n    $contentRange = $d.Range($p2.Range.Start, $p2.Range.End - 1)\n    if ($contentRange.End -gt $contentRange.Start) {\n        $contentRange.Text = \"\"\n    }\n\n    # Delete paragraph 1's ending paragraph mark to join it with paragraph 2\n    # (now empty) - this is the standard \"merge two paragraphs\" move and\n    # leaves the bookmark anchored right after the joined text.\n    $p1Again = $d.Paragraphs(1)\n    $markRange = $d.Range($p1Again.Range.End - 1, $p1Again.Range.End)\n    $markRange.Delete()\n}\n"}
